$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 8, 13, 14: touched but left blank (no formatting change) ----
$ws.Range("A8:E8").Borders.LineStyle = -4142
$ws.Range("A13:E14").Borders.LineStyle = -4142

# ---- Apply shared formatting per style group (border + alignment + font) ----
# style group s4
$r = $ws.Range("A9:A12")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("B9")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("C9")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("D9:D12")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("E9:E11")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

# style group s5
$r = $ws.Range("A16")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Color = 32768
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("A19")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Color = 32768
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("A21:A23")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Color = 32768
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("A26:A33")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Color = 32768
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("A37:A39")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Color = 32768
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("B10:B12")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Color = 32768
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("D16:D17")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Color = 32768
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

# style group s6
$r = $ws.Range("A17:A18")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Color = 255
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("A20")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Color = 255
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("A24:A25")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Color = 255
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("A34:A36")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Color = 255
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("A40")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Color = 255
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("C10:C12")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Color = 255
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("D18")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Color = 255
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

# style group s7
$r = $ws.Range("B16:B40")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Color = 16711680
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("E12")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Color = 16711680
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("E16:E18")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Color = 16711680
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

# style group s8
$r = $ws.Range("A15")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Bold = $true
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("B15")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Bold = $true
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("D15")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Bold = $true
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

$r = $ws.Range("E15")
$r.Font.Name = "Century"
$r.Font.Size = 12
$r.Font.Bold = $true
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108

# ---- Cell values ----
$ws.Range("B9").Value = "Right"
$ws.Range("C9").Value = "Wrong"
$ws.Range("D9").Value = "Not Attempt"
$ws.Range("E9").Value = "Max"
$ws.Range("A10").Value = "No."
$ws.Range("B10").Value = 18
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 8
$ws.Range("E10").Value = 28
$ws.Range("A11").Value = "Marking"
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0
$ws.Range("A12").Value = "Total"
$ws.Range("B12").Value = 90
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "88/140"
$ws.Range("A15").Value = "Student Ans"
$ws.Range("B15").Value = "Correct Ans"
$ws.Range("D15").Value = "Student Ans"
$ws.Range("E15").Value = "Correct Ans"
$ws.Range("A16").Value = "Option A"
$ws.Range("B16").Value = "Option A"
$ws.Range("D16").Value = "Option A"
$ws.Range("E16").Value = "Option A"
$ws.Range("B17").Value = "Option D"
$ws.Range("D17").Value = "Option C"
$ws.Range("E17").Value = "Option C"
$ws.Range("B18").Value = "Option B"
$ws.Range("D18").Value = "Option B"
$ws.Range("E18").Value = "Option D"
$ws.Range("A19").Value = "Option C"
$ws.Range("B19").Value = "Option C"
$ws.Range("B20").Value = "Option B"
$ws.Range("A21").Value = "Option C"
$ws.Range("B21").Value = "Option C"
$ws.Range("A22").Value = "Option D"
$ws.Range("B22").Value = "Option D"
$ws.Range("A23").Value = "Option D"
$ws.Range("B23").Value = "Option D"
$ws.Range("B24").Value = "Option A"
$ws.Range("B25").Value = "Option A"
$ws.Range("A26").Value = "Option C"
$ws.Range("B26").Value = "Option C"
$ws.Range("A27").Value = "Option A"
$ws.Range("B27").Value = "Option A"
$ws.Range("A28").Value = "Option D"
$ws.Range("B28").Value = "Option D"
$ws.Range("A29").Value = "Option D"
$ws.Range("B29").Value = "Option D"
$ws.Range("A30").Value = "Option B"
$ws.Range("B30").Value = "Option B"
$ws.Range("A31").Value = "Option D"
$ws.Range("B31").Value = "Option D"
$ws.Range("A32").Value = "Option C"
$ws.Range("B32").Value = "Option C"
$ws.Range("A33").Value = "Option D"
$ws.Range("B33").Value = "Option D"
$ws.Range("B34").Value = "Option B"
$ws.Range("B35").Value = "Option D"
$ws.Range("B36").Value = "Option A"
$ws.Range("A37").Value = "Option A"
$ws.Range("B37").Value = "Option A"
$ws.Range("A38").Value = "Option A"
$ws.Range("B38").Value = "Option A"
$ws.Range("A39").Value = "Option D"
$ws.Range("B39").Value = "Option D"
$ws.Range("A40").Value = "Option B"
$ws.Range("B40").Value = "Option D"
